# -----------------------------------------------------------------------
# Pituitary hormone deficiency workbook update
#  1. Add a new "metadata" worksheet after the existing "data" sheet
#  2. Populate the metadata sheet with panel-level metadata
#  3. Refresh the "time_taken" timestamps (column F) on the "data" sheet
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# --- 1. Refresh the per-row "time_taken" timestamps on the data sheet ---
$ws1.Range("F2").Value = "2021-10-05 14:35:15.513360"
$ws1.Range("F3").Value = "2021-10-05 14:35:15.513368"
$ws1.Range("F4").Value = "2021-10-05 14:35:15.513371"
$ws1.Range("F5").Value = "2021-10-05 14:35:15.513374"
$ws1.Range("F6").Value = "2021-10-05 14:35:15.513376"
$ws1.Range("F7").Value = "2021-10-05 14:35:15.513379"
$ws1.Range("F8").Value = "2021-10-05 14:35:15.513381"
$ws1.Range("F9").Value = "2021-10-05 14:35:15.513384"
$ws1.Range("F10").Value = "2021-10-05 14:35:15.513387"
$ws1.Range("F11").Value = "2021-10-05 14:35:15.513389"
$ws1.Range("F12").Value = "2021-10-05 14:35:15.513392"
$ws1.Range("F13").Value = "2021-10-05 14:35:15.513394"
$ws1.Range("F14").Value = "2021-10-05 14:35:15.513396"
$ws1.Range("F15").Value = "2021-10-05 14:35:15.513399"
$ws1.Range("F16").Value = "2021-10-05 14:35:15.513401"
$ws1.Range("F17").Value = "2021-10-05 14:35:15.513404"
$ws1.Range("F18").Value = "2021-10-05 14:35:15.513407"
$ws1.Range("F19").Value = "2021-10-05 14:35:15.513409"
$ws1.Range("F20").Value = "2021-10-05 14:35:15.513411"
$ws1.Range("F21").Value = "2021-10-05 14:35:15.513414"
$ws1.Range("F22").Value = "2021-10-05 14:35:15.513416"
$ws1.Range("F23").Value = "2021-10-05 14:35:15.513419"
$ws1.Range("F24").Value = "2021-10-05 14:35:15.513421"
$ws1.Range("F25").Value = "2021-10-05 14:35:15.513424"
$ws1.Range("F26").Value = "2021-10-05 14:35:15.513426"
$ws1.Range("F27").Value = "2021-10-05 14:35:15.513429"
$ws1.Range("F28").Value = "2021-10-05 14:35:15.513432"
$ws1.Range("F29").Value = "2021-10-05 14:35:15.513434"
$ws1.Range("F30").Value = "2021-10-05 14:35:15.513436"
$ws1.Range("F31").Value = "2021-10-05 14:35:15.513439"
$ws1.Range("F32").Value = "2021-10-05 14:35:15.513441"
$ws1.Range("F33").Value = "2021-10-05 14:35:15.513444"
$ws1.Range("F34").Value = "2021-10-05 14:35:15.513447"
$ws1.Range("F35").Value = "2021-10-05 14:35:15.513449"
$ws1.Range("F36").Value = "2021-10-05 14:35:15.513452"
$ws1.Range("F37").Value = "2021-10-05 14:35:15.513454"
$ws1.Range("F38").Value = "2021-10-05 14:35:15.513457"
$ws1.Range("F39").Value = "2021-10-05 14:35:15.513459"
$ws1.Range("F40").Value = "2021-10-05 14:35:15.513462"
$ws1.Range("F41").Value = "2021-10-05 14:35:15.513464"
$ws1.Range("F42").Value = "2021-10-05 14:35:15.513467"
$ws1.Range("F43").Value = "2021-10-05 14:35:15.513470"
$ws1.Range("F44").Value = "2021-10-05 14:35:15.513472"
$ws1.Range("F45").Value = "2021-10-05 14:35:15.513474"
$ws1.Range("F46").Value = "2021-10-05 14:35:15.513477"
$ws1.Range("F47").Value = "2021-10-05 14:35:15.513479"
$ws1.Range("F48").Value = "2021-10-05 14:35:15.513482"
$ws1.Range("F49").Value = "2021-10-05 14:35:15.513484"
$ws1.Range("F50").Value = "2021-10-05 14:35:15.513487"
$ws1.Range("F51").Value = "2021-10-05 14:35:15.513489"
$ws1.Range("F52").Value = "2021-10-05 14:35:15.513491"
$ws1.Range("F53").Value = "2021-10-05 14:35:15.513494"
$ws1.Range("F54").Value = "2021-10-05 14:35:15.513497"
$ws1.Range("F55").Value = "2021-10-05 14:35:15.513499"
$ws1.Range("F56").Value = "2021-10-05 14:35:15.513502"

# --- 2. Add the new "metadata" worksheet right after "data" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Pituitary hormone deficiency"
$ws2.Range("C2").Value = 3236
$ws2.Range("D2").Value = "'0.20"
$ws2.Range("E2").Value = "2021-10-01T08:28:44.234979Z"
$ws2.Range("F2").Value = "2021-10-05 14:35:15.509615"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3236/?format=json"

# Match the bold/bordered header style used on the "data" sheet's header row,
# and the style used for the leading index column (A2), by copying formats.
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# --- 3. Keep "data" as the active sheet, as in the original workbook ---
$ws1.Activate()
$ws1.Range("A1").Select()
